# Apply weekly-cut update: revenue figures were corrected from $9,591.40 to $8.06
# (and the matching top-product quantity/revenue row updated accordingly).

$wb = $excel.ActiveWorkbook

# --- Sheet "Resumen": update the two textual revenue cells ---
# Force the cells to be treated as plain text (instead of letting Excel
# auto-convert the "$"-prefixed value into a currency-formatted number),
# then restore the "Normal" style so no extra formatting is left behind.
$wsResumen = $wb.Worksheets.Item("Resumen")

$wsResumen.Range("B3").NumberFormat = "@"
$wsResumen.Range("B3").Value = "$8.06"
$wsResumen.Range("B3").Style = "Normal"

$wsResumen.Range("B6").NumberFormat = "@"
$wsResumen.Range("B6").Value = "$8.06"
$wsResumen.Range("B6").Style = "Normal"

# --- Sheet "Top Productos": update quantity and revenue for the top product row ---
$wsTop = $wb.Worksheets.Item("Top Productos")
$wsTop.Range("B2").Value = 1
$wsTop.Range("C2").Value = 8.06

Write-Output "Updated Resumen!B3, Resumen!B6, Top Productos!B2, Top Productos!C2"
